# Kalgor Advancements fix:
# "fixed kalgor advancements to apply from prev advancements to future weapons"
#
# - blade of dread (H/I) damage now keeps increasing for strength1-3 rows
# - tenebris touch (K/L) now also applies to strength1-3 rows (6-8), and row9 value increases
# - tenebris touch no longer applies to row10 (abyss2)
# - ring of darkness (N/O) now also applies to strength1-3 rows (6-8) and to row10 (abyss2)
# - ring of darkness no longer applies to rows 11/12 (which get removed / merged away)
# - "abyss1"/"abyss2" collapse into a single "abyss"/"darkness" label (rows 9-10); row11/12 names removed
# - unused trailing rows 21 and 22 are removed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- strength1 (row 6): blade of dread damage bump + new tenebris touch / ring of darkness values ---
$ws.Range("H6").Value2 = 15
$ws.Range("K6").Value2 = 11
$ws.Range("L6").Value2 = 2
$ws.Range("N6").Value2 = 7
$ws.Range("O6").Value2 = 3

# --- strength2 (row 7): blade of dread damage bump + new tenebris touch / ring of darkness values ---
$ws.Range("H7").Value2 = 17
$ws.Range("K7").Value2 = 12
$ws.Range("L7").Value2 = 2
$ws.Range("N7").Value2 = 9
$ws.Range("O7").Value2 = 3

# --- strength3 (row 8): blade of dread damage bump + new tenebris touch / ring of darkness values ---
$ws.Range("H8").Value2 = 19
$ws.Range("K8").Value2 = 13
$ws.Range("L8").Value2 = 2
$ws.Range("N8").Value2 = 11
$ws.Range("O8").Value2 = 3

# --- row 9 (abyss1 -> abyss): tenebris touch increases ---
$ws.Range("A9").Value2 = "abyss"
$ws.Range("K9").Value2 = 13

# --- row 10 (abyss2 -> darkness): tenebris touch removed, ring of darkness added instead ---
$ws.Range("A10").Value2 = "darkness"
$ws.Range("K10:L10").ClearContents() | Out-Null
$ws.Range("N10").Value2 = 12
$ws.Range("O10").Value2 = 4

# --- row 11 (was darkness1): label and ring of darkness values removed ---
$ws.Range("A11").ClearContents() | Out-Null
$ws.Range("N11:O11").ClearContents() | Out-Null

# --- row 12 (was darkness2): label and ring of darkness values removed ---
$ws.Range("A12").ClearContents() | Out-Null
$ws.Range("N12:O12").ClearContents() | Out-Null

# --- remove now-unused trailing rows 21 and 22 ---
$ws.Rows("21:22").Delete() | Out-Null

# --- update selection to match saved cursor position ---
$ws.Range("A11").Select() | Out-Null
